$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 5774
$ws1.Range("F8").Value  = 6376
$ws1.Range("F10").Value = 1942
$ws1.Range("F18").Value = 7979
$ws1.Range("F19").Value = 7979
$ws1.Range("F20").Value = 140
$ws1.Range("F31").Value = 1787

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 205

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 5774
$ws4.Range("F14").Value = 6376
$ws4.Range("F16").Value = 1942
$ws4.Range("F24").Value = 7979
$ws4.Range("F25").Value = 7979
$ws4.Range("F26").Value = 140
$ws4.Range("F36").Value = 1787
